$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarterly-label values that replace the raw date serials in column A
$labels = @(
    "2004Q4",
    "2005Q4",
    "2006Q4",
    "2007Q4",
    "2008Q4",
    "2009Q4",
    "2010Q4",
    "2011Q4",
    "2012Q4",
    "2013Q4",
    "2014Q4",
    "2015Q4",
    "2016Q4",
    "2017Q4",
    "2018Q4",
    "2019Q4",
    "2020Q4",
    "2021Q4",
    "2022Q4",
    "2023Q4",
    "2024Q4"
)

# Copy the header cell's formatting (centered, bordered, bold, General
# number format) onto the data cells in column A so they no longer use
# the custom date/time number format.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Replace each date serial with its quarter-label text equivalent.
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}
